$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $val)
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "56.876.16"
Set-TextValue $ws "E2" "  -0.78%  "
Set-TextValue $ws "D3" "2.971.74"
Set-TextValue $ws "E3" "  -1.82%  "
Set-TextValue $ws "E4" "  +0.13%  "
Set-TextValue $ws "D5" "499.31"
Set-TextValue $ws "E5" "  -3.92%  "
Set-TextValue $ws "D6" "137.07"
Set-TextValue $ws "E6" "  -3.10%  "
Set-TextValue $ws "E7" "  +0.09%  "
Set-TextValue $ws "D8" "0.428"
Set-TextValue $ws "E8" "  -2.36%  "
Set-TextValue $ws "D9" "7.31"
Set-TextValue $ws "E9" "  -3.77%  "
Set-TextValue $ws "E10" "  -2.53%  "
Set-TextValue $ws "E11" "  -1.14%  "
Set-TextValue $ws "D12" "3.477.38"
Set-TextValue $ws "E12" "  -1.88%  "
Set-TextValue $ws "D13" "0.127"
Set-TextValue $ws "E13" "  -1.97%  "
Set-TextValue $ws "D14" "25.87"
Set-TextValue $ws "E14" "  -1.35%  "
Set-TextValue $ws "E15" "  -1.70%  "
Set-TextValue $ws "D16" "56.973.17"
Set-TextValue $ws "E16" "  -0.57%  "
Set-TextValue $ws "D17" "6.04"
Set-TextValue $ws "E17" "  +0.10%  "
Set-TextValue $ws "D18" "2.978.51"
Set-TextValue $ws "E18" "  -1.37%  "
Set-TextValue $ws "D19" "12.55"
Set-TextValue $ws "E19" "  -1.46%  "
Set-TextValue $ws "D20" "7.81"
Set-TextValue $ws "E20" "  -2.27%  "
Set-TextValue $ws "D21" "318.65"
Set-TextValue $ws "E21" "  -4.31%  "
Set-TextValue $ws "E22" "  -0.13%  "
Set-TextValue $ws "E23" "  -0.83%  "
Set-TextValue $ws "D24" "0.485"
Set-TextValue $ws "E24" "  -0.88%  "
Set-TextValue $ws "D25" "63.32"
Set-TextValue $ws "E25" "  -1.26%  "
Set-TextValue $ws "E26" "  +0.24%  "
Set-TextValue $ws "E27" "  -5.03%  "
Set-TextValue $ws "E28" "  -4.20%  "
Set-TextValue $ws "E29" "  -4.37%  "
Set-TextValue $ws "E30" "  -2.42%  "
Set-TextValue $ws "E31" "  -3.69%  "
Set-TextValue $ws "D32" "1.15"
Set-TextValue $ws "E32" "  -6.19%  "
Set-TextValue $ws "D33" "20.11"
Set-TextValue $ws "E33" "  -3.66%  "
Set-TextValue $ws "D34" "156.41"
Set-TextValue $ws "E34" "  -1.42%  "
Set-TextValue $ws "D35" "4.56"
Set-TextValue $ws "E35" "  -2.33%  "
Set-TextValue $ws "D36" "5.74"
Set-TextValue $ws "E36" "  -1.02%  "
Set-TextValue $ws "E37" "  -5.16%  "
Set-TextValue $ws "D38" "24.04"
Set-TextValue $ws "E38" "  -1.94%  "
Set-TextValue $ws "D39" "0.0663"
Set-TextValue $ws "E39" "  -2.45%  "
Set-TextValue $ws "D40" "3.004.21"
Set-TextValue $ws "E40" "  -1.76%  "
Set-TextValue $ws "D41" "37.52"
Set-TextValue $ws "E41" "  +0.12%  "
Set-TextValue $ws "E42" "  -0.02%  "
Set-TextValue $ws "D43" "3.72"
Set-TextValue $ws "E43" "  -1.10%  "
Set-TextValue $ws "D44" "0.637"
Set-TextValue $ws "E44" "  -2.71%  "
Set-TextValue $ws "D45" "2.201.06"
Set-TextValue $ws "E45" "  -4.46%  "
Set-TextValue $ws "E46" "  -4.02%  "
Set-TextValue $ws "B47" "Cosmos"
Set-TextValue $ws "C47" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws "D47" "5.92"
Set-TextValue $ws "E47" "  -0.21%  "
Set-TextValue $ws "B48" "ONDO"
Set-TextValue $ws "C48" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws "D48" "0.936"
Set-TextValue $ws "E48" "  -8.00%  "
Set-TextValue $ws "E49" "  -3.65%  "
Set-TextValue $ws "D50" "19.13"
Set-TextValue $ws "E50" "  -2.26%  "
Set-TextValue $ws "D51" "1.80"
Set-TextValue $ws "E51" "  -11.33%  "
